$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.543.82'
$ws.Range("E2").Value = '  -2.36%  '

$ws.Range("D3").Value = '3.315.27'
$ws.Range("E3").Value = '  -4.18%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '547.68'
$ws.Range("E5").Value = '  -0.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.89'
$ws.Range("E6").Value = '  -3.62%  '

$ws.Range("E7").Value = '  -4.13%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("E9").Value = '  -3.26%  '

$ws.Range("E10").Value = '  -0.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.07'
$ws.Range("E11").Value = '  -1.10%  '

$ws.Range("E12").Value = '  -1.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.83'
$ws.Range("E13").Value = '  -3.96%  '

$ws.Range("D14").Value = '3.851.63'
$ws.Range("E14").Value = '  -4.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.03'
$ws.Range("E15").Value = '  -2.94%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.323.18'
$ws.Range("E16").Value = '  -4.33%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.117'
$ws.Range("E17").Value = '  -3.61%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '63.544.14'
$ws.Range("E18").Value = '  -2.60%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.58'
$ws.Range("E19").Value = '  -2.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.971'
$ws.Range("E20").Value = '  -1.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '411.55'
$ws.Range("E21").Value = '  -0.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.02'
$ws.Range("E22").Value = '  -0.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.36'
$ws.Range("E23").Value = '  +6.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.64'
$ws.Range("E24").Value = '  +6.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.66'
$ws.Range("E25").Value = '  -3.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.50'
$ws.Range("E26").Value = '  -2.72%  '

$ws.Range("E27").Value = '  -4.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.57'
$ws.Range("E28").Value = '  -5.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '28.95'
$ws.Range("E29").Value = '  -3.96%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.35'
$ws.Range("E30").Value = '  -2.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.30'
$ws.Range("E31").Value = '  -3.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '573.95'
$ws.Range("E32").Value = '  -6.35%  '

$ws.Range("E33").Value = '  -3.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.57'
$ws.Range("E34").Value = '  -2.28%  '

$ws.Range("E35").Value = '  +0.07%  '

$ws.Range("E36").Value = '  +1.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '34.97'
$ws.Range("E37").Value = '  -6.31%  '

$ws.Range("E38").Value = '  +4.26%  '

$ws.Range("E39").Value = '  -6.37%  '

$ws.Range("E40").Value = '  -3.89%  '

$ws.Range("D41").Value = '3.112.92'
$ws.Range("E41").Value = '  -7.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.26'
$ws.Range("E43").Value = '  +1.21%  '

$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.77'
$ws.Range("E44").Value = '  -1.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0398'
$ws.Range("E45").Value = '  -3.14%  '

$ws.Range("E46").Value = '  -4.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.60'
$ws.Range("E47").Value = '  -4.27%  '

$ws.Range("E48").Value = '  -3.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.70'
$ws.Range("E49").Value = '  -3.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.97'
$ws.Range("E50").Value = '  -5.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000226'
$ws.Range("E51").Value = '  +5.13%  '
